$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1086
$ws.Range("J17").Value = 1086
$ws.Range("L17").Value = 3258
$ws.Range("N17").Value = -3594

$ws.Range("H62").Value = 90950490
$ws.Range("I62").Value = 250001250
$ws.Range("K62").Value = 250001250
$ws.Range("M62").Value = -250000626

$ws.Range("H65").Value = 90950490
$ws.Range("I65").Value = 250001250
$ws.Range("K65").Value = 1250006250
$ws.Range("M65").Value = -1250003130

$ws.Range("H80").Value = 39451.69
$ws.Range("I80").Value = 16828.834
$ws.Range("J80").Value = 58842.715
$ws.Range("K80").Value = 50486.50199999999
$ws.Range("L80").Value = 176528.145
$ws.Range("M80").Value = -49488.50199999999
$ws.Range("N80").Value = -178524.145

$ws.Range("H83").Value = 39451.69
$ws.Range("I83").Value = 16828.834
$ws.Range("J83").Value = 58842.715
$ws.Range("K83").Value = 151459.506
$ws.Range("L83").Value = 529584.4349999999
$ws.Range("M83").Value = -146467.506
$ws.Range("N83").Value = -539568.4349999999

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H129").Value = 1769
$ws.Range("I129").Value = 1769
$ws.Range("K129").Value = 5307
$ws.Range("M129").Value = -307

$ws.Range("H132").Value = 1872.1451
$ws.Range("I132").Value = 1831.7627
$ws.Range("K132").Value = 5495.2881
$ws.Range("M132").Value = -2965.2881

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2029.65
$ws.Range("I2").Value = 1792.1428
$ws.Range("J2").Value = 2583.8333
$ws.Range("K2").Value = 1792.1428
$ws.Range("L2").Value = 2583.8333
$ws.Range("M2").Value = -1679.1428
$ws.Range("N2").Value = -2809.8333

$ws.Range("H32").Value = 1669397.4
$ws.Range("I32").Value = 1814214.6
$ws.Range("K32").Value = 1814214.6
$ws.Range("M32").Value = -1813927.6

$ws.Range("H39").Value = 1508
$ws.Range("I39").Value = 1508
$ws.Range("K39").Value = 1508
$ws.Range("M39").Value = -988

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H97").Value = 4637051
$ws.Range("I97").Value = 818.1667
$ws.Range("K97").Value = 818.1667
$ws.Range("M97").Value = -322.1667

$ws.Range("H102").Value = 4636.75
$ws.Range("I102").Value = 4449.3335
$ws.Range("K102").Value = 4449.3335
$ws.Range("M102").Value = -2827.3335

$ws.Range("H116").Value = 2029.65
$ws.Range("I116").Value = 1792.1428
$ws.Range("J116").Value = 2583.8333
$ws.Range("K116").Value = 1792.1428
$ws.Range("L116").Value = 2583.8333
$ws.Range("M116").Value = 501.8571999999999
$ws.Range("N116").Value = -7171.8333

$ws.Range("H122").Value = 23326.818
$ws.Range("I122").Value = 29075
$ws.Range("K122").Value = 87225
$ws.Range("M122").Value = -84775

$ws.Range("H132").Value = 8509.727999999999
$ws.Range("I132").Value = 3665.2856
$ws.Range("J132").Value = 16987.5
$ws.Range("K132").Value = 10995.8568
$ws.Range("L132").Value = 50962.5
$ws.Range("M132").Value = -8465.856800000001
$ws.Range("N132").Value = -56022.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2029.65
$ws.Range("I3").Value = 1792.1428
$ws.Range("J3").Value = 2583.8333
$ws.Range("K3").Value = 1792.1428
$ws.Range("L3").Value = 2583.8333
$ws.Range("M3").Value = -1678.1428
$ws.Range("N3").Value = -2811.8333

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H94").Value = 1469.8334
$ws.Range("I94").Value = 763.8
$ws.Range("K94").Value = 763.8
$ws.Range("M94").Value = -312.8

$ws.Range("H105").Value = 2865.3076
$ws.Range("I105").Value = 1933.2
$ws.Range("K105").Value = 1933.2
$ws.Range("M105").Value = -186.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 10000
$ws.Range("I39").Value = 10000
$ws.Range("K39").Value = 10000
$ws.Range("M39").Value = -9609

$ws.Range("H49").Value = 10000
$ws.Range("I49").Value = 10000
$ws.Range("K49").Value = 10000
$ws.Range("M49").Value = -9818

$ws.Range("H62").Value = 6275.8
$ws.Range("J62").Value = 9000
$ws.Range("L62").Value = 9000
$ws.Range("N62").Value = -10248

$ws.Range("H65").Value = 6275.8
$ws.Range("J65").Value = 9000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51240

$ws.Range("H132").Value = 9611.125
$ws.Range("I132").Value = 2697
$ws.Range("J132").Value = 11915.833
$ws.Range("K132").Value = 8091
$ws.Range("L132").Value = 35747.499
$ws.Range("M132").Value = -5561
$ws.Range("N132").Value = -40807.499

$ws.Range("H138").Value = 99999.5
$ws.Range("J138").Value = 99999.5
$ws.Range("L138").Value = 99999.5
$ws.Range("N138").Value = -110279.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 13994.5
$ws.Range("J132").Value = 25249.75
$ws.Range("L132").Value = 227247.75
$ws.Range("N132").Value = -232307.75

$ws.Range("H134").Value = 83663.84
$ws.Range("I134").Value = 106263
$ws.Range("J134").Value = 8333.333000000001
$ws.Range("K134").Value = 318789
$ws.Range("L134").Value = 24999.999
$ws.Range("M134").Value = -313719
$ws.Range("N134").Value = -35139.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5813
$ws.Range("I132").Value = 2401.1
$ws.Range("J132").Value = 11499.5
$ws.Range("K132").Value = 7203.299999999999
$ws.Range("L132").Value = 34498.5
$ws.Range("M132").Value = -4673.299999999999
$ws.Range("N132").Value = -39558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 4030
$ws.Range("I14").Value = 90
$ws.Range("K14").Value = 90
$ws.Range("M14").Value = 82

$ws.Range("H16").Value = 1023.9655
$ws.Range("I16").Value = 1107.04
$ws.Range("J16").Value = 504.75
$ws.Range("K16").Value = 1107.04
$ws.Range("L16").Value = 504.75
$ws.Range("M16").Value = -937.04
$ws.Range("N16").Value = -844.75

$ws.Range("H122").Value = 3797.9707
$ws.Range("I122").Value = 2875
$ws.Range("K122").Value = 8625
$ws.Range("M122").Value = -6175

$ws.Range("H132").Value = 17865808
$ws.Range("I132").Value = 38466660
$ws.Range("J132").Value = 11733.333
$ws.Range("K132").Value = 115399980
$ws.Range("L132").Value = 35199.999
$ws.Range("M132").Value = -115397450
$ws.Range("N132").Value = -40259.999

$ws.Range("H136").Value = 10321.107
$ws.Range("I136").Value = 5289
$ws.Range("J136").Value = 12704.737
$ws.Range("K136").Value = 15867
$ws.Range("L136").Value = 38114.211
$ws.Range("M136").Value = -13317
$ws.Range("N136").Value = -43214.211

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 500000670
$ws.Range("I14").Value = 500000670
$ws.Range("K14").Value = 500000670
$ws.Range("M14").Value = -500000502

$ws.Range("H21").Value = 20000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20470

$ws.Range("H33").Value = 8000
$ws.Range("I33").Value = 6000
$ws.Range("K33").Value = 6000
$ws.Range("M33").Value = -5750

$ws.Range("H35").Value = 20000
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20580

$ws.Range("H36").Value = 8000
$ws.Range("I36").Value = 6000
$ws.Range("K36").Value = 6000
$ws.Range("M36").Value = -5750

$ws.Range("H95").Value = 400000
$ws.Range("J95").Value = 400000
$ws.Range("L95").Value = 400000
$ws.Range("N95").Value = -405492

$ws.Range("H107").Value = 851.5
$ws.Range("I107").Value = 635.44446
$ws.Range("K107").Value = 1906.33338
$ws.Range("M107").Value = 13.66661999999997

$ws.Range("H122").Value = 22912854
$ws.Range("I122").Value = 31502392
$ws.Range("K122").Value = 94507176
$ws.Range("M122").Value = -94504726

$ws.Range("H132").Value = 100120500
$ws.Range("I132").Value = 250050000
$ws.Range("J132").Value = 167499.67
$ws.Range("K132").Value = 750150000
$ws.Range("L132").Value = 502499.01
$ws.Range("M132").Value = -750147470
$ws.Range("N132").Value = -507559.01

$ws.Range("H136").Value = 40043770
$ws.Range("I136").Value = 55557490
$ws.Range("K136").Value = 166672470
$ws.Range("M136").Value = -166669920
